# Update "想去人数" (number of people wanting to go) counts on the
# "展览" and "全部类型" sheets to match the regenerated site data.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws4 = $wb.Worksheets.Item("全部类型")

# "展览" sheet (sheet1)
$ws1.Range("F3").Value = 2204
$ws1.Range("F4").Value = 88
$ws1.Range("F5").Value = 13171
$ws1.Range("F8").Value = 517
$ws1.Range("F10").Value = 1184
$ws1.Range("F11").Value = 989
$ws1.Range("F12").Value = 13775
$ws1.Range("F13").Value = 14384
$ws1.Range("F23").Value = 112
$ws1.Range("F25").Value = 5438
$ws1.Range("F26").Value = 940
$ws1.Range("F27").Value = 23
$ws1.Range("F28").Value = 321
$ws1.Range("F29").Value = 23
$ws1.Range("F30").Value = 57

# "全部类型" sheet (sheet4)
$ws4.Range("F3").Value = 2204
$ws4.Range("F4").Value = 88
$ws4.Range("F5").Value = 13172
$ws4.Range("F8").Value = 517
$ws4.Range("F10").Value = 1184
$ws4.Range("F11").Value = 989
$ws4.Range("F12").Value = 13775
$ws4.Range("F13").Value = 14384
$ws4.Range("F23").Value = 112
$ws4.Range("F25").Value = 5438
$ws4.Range("F26").Value = 940
$ws4.Range("F27").Value = 23
$ws4.Range("F28").Value = 321
$ws4.Range("F29").Value = 23
$ws4.Range("F30").Value = 57
